# Fruta / hortaliza, semanal
# Insert a new weekly observation row at sheet row 222 (pushing the
# existing rows 222-288 down to 223-289) and populate it with the new
# data point (Fukumoto / Primera, week of 2022-06-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 222; this shifts every
# row from 222 down to 288 one position down (to 223-289) and extends
# the used range / dimension accordingly.
$ws.Rows.Item(222).Insert()

$ws.Cells.Item(222, 1).Value = 11
$ws.Cells.Item(222, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(222, 3).Value = "Bíobío"
$ws.Cells.Item(222, 4).Value = 44726
$ws.Cells.Item(222, 5).Value = 8
$ws.Cells.Item(222, 6).Value = "Fruta"
$ws.Cells.Item(222, 7).Value = 100102
$ws.Cells.Item(222, 8).Value = "Cítricos"
$ws.Cells.Item(222, 9).Value = 100102005
$ws.Cells.Item(222, 10).Value = "Naranja"
$ws.Cells.Item(222, 11).Value = "Fukumoto"
$ws.Cells.Item(222, 12).Value = "Primera"
$ws.Cells.Item(222, 13).Value = 100
$ws.Cells.Item(222, 14).Value = 7000
$ws.Cells.Item(222, 15).Value = 8000
$ws.Cells.Item(222, 16).Value = 7500
$ws.Cells.Item(222, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(222, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(222, 19).Value = 500
$ws.Cells.Item(222, 20).Value = 15
